$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (want-to-go count) column F for rows 3-5
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 149
$wsExhibit.Range("F4").Value = 704
$wsExhibit.Range("F5").Value = 62

# Sheet "全部类型" - same events appear one row lower (rows 4-6)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 149
$wsAll.Range("F5").Value = 704
$wsAll.Range("F6").Value = 62
